$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextDate($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextDate "A3" "28-07-2022"
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

Set-TextDate "A4" "01-08-2022"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

Set-TextDate "A5" "04-08-2022"
Set-TextDate "A6" "08-08-2022"
Set-TextDate "A7" "11-08-2022"
Set-TextDate "A8" "15-08-2022"
Set-TextDate "A9" "18-08-2022"
Set-TextDate "A10" "22-08-2022"
Set-TextDate "A11" "25-08-2022"
Set-TextDate "A12" "29-08-2022"
Set-TextDate "A13" "01-09-2022"
Set-TextDate "A14" "05-09-2022"
Set-TextDate "A15" "08-09-2022"
Set-TextDate "A16" "12-09-2022"
Set-TextDate "A17" "15-09-2022"
Set-TextDate "A18" "19-09-2022"
Set-TextDate "A19" "22-09-2022"
Set-TextDate "A20" "26-09-2022"
Set-TextDate "A21" "29-09-2022"
